# Insert a new row at position 620, shifting existing rows 620:686 down to 621:687.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(620).Insert()

# Populate the newly inserted row 620 with the new record's data.
$ws.Range("A620").Value = 5
$ws.Range("B620").Value = "Macroferia Regional de Talca"
$ws.Range("C620").Value = "Maule"
$ws.Range("D620").Value = 44946
$ws.Range("E620").Value = 7
$ws.Range("F620").Value = 100114001
$ws.Range("G620").Value = "Papa"
$ws.Range("H620").Value = "Rodeo"
$ws.Range("I620").Value = "1a (cosecha)"
$ws.Range("J620").Value = 700
$ws.Range("K620").Value = 9500
$ws.Range("L620").Value = 10000
$ws.Range("M620").Value = 9786
$ws.Range("N620").Value = "`$/saco 25 kilos"
$ws.Range("O620").Value = "Región del Maule"
$ws.Range("P620").Value = 391
$ws.Range("Q620").Value = 25
$ws.Range("R620").Value = "Hortaliza"
